$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '27.045.44'
$ws.Range("D2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '1.821.46'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +3.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = '1.010'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.88%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '314.64'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.26%  '

$ws.Range("E6").Value = '  +0.96%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.4316'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '0.3702'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.93%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.07286'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.91%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '2.140.38'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +22.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '0.8726'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '21.35'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +5.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '6.658'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '5.423'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +3.21%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '0.06972'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '81.26'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '1.012'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.83%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '0.000008872'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '1.010'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.86%  '

$ws.Range("E20").Value = '  +1.71%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '27.119.30'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.96%  '

$ws.Range("E22").Value = '  +4.06%  '

$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '11.04'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.57%  '

$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '2.386.28'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +21.46%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '154.92'
$ws.Range("D25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '1.894'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '18.46'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.56%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '5.256'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '1.924'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +11.88%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '115.09'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.66%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '0.08989'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '1.183'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +6.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '0.7498'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.71%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '4.442'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '2.826'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.01%  '

$ws.Range("E36").Value = '  +1.06%  '

$ws.Range("E37").Value = '  +5.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '0.05254'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.49%  '

$ws.Range("E39").Value = '  +2.26%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '0.5135'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '2.757'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +9.69%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '0.1657'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '6.502'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '8.334'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.44%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '107.44'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '10.46'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.76%  '

$ws.Range("E47").Value = '  +1.12%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '1.659'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +4.92%  '

$ws.Range("B49").Value = 'Decentraland'
$ws.Range("C49").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '0.4579'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.86%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '0.06241'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.77%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '1.834'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.48%  '
